$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.049.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -2.21%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.826.13'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.07%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.87%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''311.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -2.14%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = '''  -0.87%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.4236'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -1.44%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -2.03%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.07223'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.38%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.8445'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -3.39%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''20.73'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -3.56%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''1.824.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -1.26%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''6.657'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.11%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''5.294'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.48%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.07045'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -1.08%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''89.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.86%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.11%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''0.000008746'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.83%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''1.001'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.86%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''14.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -3.47%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''27.088.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.16%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''5.133'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.41%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  -1.87%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''2.049.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -1.26%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''1.981'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.49%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -2.46%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''2.257'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +4.29%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''18.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.12%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''5.245'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -1.86%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''116.78'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -1.49%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.08704'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -2.42%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -3.62%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.7369'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -5.03%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''2.902'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -0.09%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''4.422'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -2.47%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -1.03%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -3.29%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -1.82%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.05244'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -1.70%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''7.339'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +2.91%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''2.875'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.30%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.1687'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -0.34%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.5063'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -1.19%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''8.550'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -2.64%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''10.57'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -0.67%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''1.970'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +6.80%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.4724'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.66%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''105.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -1.52%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''1.001'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -1.02%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.06325'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.97%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -2.10%  '
$ws.Range("E51").Style = "Normal"
